# Updated cryptos list on Sat Jun  1 03:15:57 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.690.26"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "3.778.45"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "3.775.83"
$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.93%  "

$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "4.413.71"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").Value = "3.782.08"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").Value = "67.669.82"
$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "456.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.76%  "

$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("E24").Value = "  +3.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.23"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "

$ws.Range("E27").Value = "  -2.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").Value = "3.733.16"
$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.98%  "

$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "46.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.82%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.298"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "388.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "

$ws.Range("E50").Value = "  -4.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.83%  "
